$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Chandra, Gremlin Wrangler', ['{2}{R}{R}', 'Legendary Planeswalker — Chandra', '+1: Create a 2/2 red Gremlin creature token.', '−2: Chandra, Gremlin Wrangler deals X damage to target creature or player, where X is the number of Gremlins you control.', 'Loyalty: 3'])"
$ws.Range("A3").Value = "('Dungeon Master', ['{2}{W}{U}', 'Legendary Planeswalker — Dungeon Master', '+1: Target opponent creates a 1/1 black Skeleton creature token with “When this creature dies, each opponent gains 2 life.”', '+1: Roll a d20. If you roll a 1, skip your next turn. If you roll a 12 or higher, draw a card.', '−6: You get an adventuring party. (Your party is a 3/3 red Fighter with first strike, a 1/1 white Cleric with lifelink, a 2/2 black Rogue with hexproof, and a 1/1 blue Wizard with flying.)', 'Loyalty: 1d4+1'])"
$ws.Range("A4").Value = "('Nira, Hellkite Duelist', ['{W}{U}{B}{R}{G}', 'Legendary Creature — Dragon', 'Flash', 'Flying, trample, haste', 'When Nira, Hellkite Duelist enters the battlefield, the next time you would lose the game this turn, instead draw three cards and your life total becomes 5.', '6/6'])"

$ws.Range("A5:A21").EntireRow.Delete()
